$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 13075
$ws.Range("I18").Value = 13075
$ws.Range("K18").Value = 13075
$ws.Range("M18").Value = -12791
$ws.Range("H33").Value = 14064.956
$ws.Range("I33").Value = 15333.048
$ws.Range("J33").Value = 750
$ws.Range("K33").Value = 15333.048
$ws.Range("L33").Value = 750
$ws.Range("M33").Value = -15104.048
$ws.Range("N33").Value = -1208
$ws.Range("H80").Value = 587
$ws.Range("I80").Value = 306.84616
$ws.Range("K80").Value = 920.5384799999999
$ws.Range("M80").Value = 77.46152000000006
$ws.Range("H83").Value = 587
$ws.Range("I83").Value = 306.84616
$ws.Range("K83").Value = 2761.61544
$ws.Range("M83").Value = 2230.38456
$ws.Range("H86").Value = 5112.5
$ws.Range("I86").Value = 5000
$ws.Range("J86").Value = 5225
$ws.Range("K86").Value = 5000
$ws.Range("L86").Value = 5225
$ws.Range("M86").Value = -3877
$ws.Range("N86").Value = -7471
$ws.Range("H89").Value = 5112.5
$ws.Range("I89").Value = 5000
$ws.Range("J89").Value = 5225
$ws.Range("K89").Value = 25000
$ws.Range("L89").Value = 26125
$ws.Range("M89").Value = -19384
$ws.Range("N89").Value = -37357
$ws.Range("H132").Value = 1908.7368
$ws.Range("I132").Value = 1924.4706
$ws.Range("J132").Value = 1775
$ws.Range("K132").Value = 5773.4118
$ws.Range("L132").Value = 5325
$ws.Range("M132").Value = -3243.4118
$ws.Range("N132").Value = -10385
$ws.Range("H137").Value = 1738.1111
$ws.Range("I137").Value = 1458.6
$ws.Range("K137").Value = 4375.799999999999
$ws.Range("M137").Value = -1825.799999999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8440.849
$ws.Range("I32").Value = 8340.258
$ws.Range("K32").Value = 8340.258
$ws.Range("M32").Value = -8053.258
$ws.Range("H61").Value = 6331.9614
$ws.Range("I61").Value = 4831.727
$ws.Range("K61").Value = 4831.727
$ws.Range("M61").Value = -4619.727
$ws.Range("H136").Value = 6331.9614
$ws.Range("I136").Value = 4831.727
$ws.Range("K136").Value = 14495.181
$ws.Range("M136").Value = -11945.181

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("N85").ClearContents()
$ws.Range("H86").Value = 1754.5385
$ws.Range("I86").Value = 1525.75
$ws.Range("K86").Value = 1525.75
$ws.Range("M86").Value = -402.75
$ws.Range("H89").Value = 1754.5385
$ws.Range("I89").Value = 1525.75
$ws.Range("K89").Value = 7628.75
$ws.Range("M89").Value = -2012.75
$ws.Range("H134").Value = 3531.9062
$ws.Range("I134").Value = 3613.5806
$ws.Range("K134").Value = 10840.7418
$ws.Range("M134").Value = -8305.7418
$ws.Range("H140").Value = 81983
$ws.Range("J140").Value = 81983
$ws.Range("L140").Value = 81983
$ws.Range("N140").Value = -92343

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 20555.428
$ws.Range("J41").Value = 20163.334
$ws.Range("L41").Value = 20163.334
$ws.Range("N41").Value = -21019.334
$ws.Range("H59").Value = 34305
$ws.Range("J59").Value = 40735
$ws.Range("L59").Value = 40735
$ws.Range("N59").Value = -43025
$ws.Range("H97").Value = 34985
$ws.Range("J97").Value = 34985
$ws.Range("L97").Value = 34985
$ws.Range("N97").Value = -36967
$ws.Range("H99").Value = 6410.4
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()
$ws.Range("H126").Value = 6410.4
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 14000
$ws.Range("J21").Value = 14000
$ws.Range("L21").Value = 14000
$ws.Range("N21").Value = -14346
$ws.Range("H30").Value = 14000
$ws.Range("J30").Value = 14000
$ws.Range("L30").Value = 14000
$ws.Range("N30").Value = -14210
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
$ws.Range("H102").Value = 5205.9
$ws.Range("I102").Value = 4451
$ws.Range("J102").Value = 12000
$ws.Range("K102").Value = 4451
$ws.Range("L102").Value = 12000
$ws.Range("M102").Value = -2829
$ws.Range("N102").Value = -15244
$ws.Range("H107").Value = 958199.1
$ws.Range("I107").Value = 1724387.1
$ws.Range("K107").Value = 1724387.1
$ws.Range("M107").Value = -1722467.1
$ws.Range("H113").Value = 5812
$ws.Range("I113").Value = 4996.6
$ws.Range("K113").Value = 4996.6
$ws.Range("M113").Value = -2826.6
$ws.Range("H122").Value = 3230.3
$ws.Range("I122").Value = 2550.6667
$ws.Range("K122").Value = 7652.000100000001
$ws.Range("M122").Value = -5202.000100000001
$ws.Range("H126").Value = 3617.4285
$ws.Range("I126").Value = 2564.4
$ws.Range("K126").Value = 7693.200000000001
$ws.Range("M126").Value = -5223.200000000001
$ws.Range("H132").Value = 3158.25
$ws.Range("I132").Value = 2517
$ws.Range("K132").Value = 7551
$ws.Range("M132").Value = -5021

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H8").Value = 130000
$ws.Range("J8").Value = 130000
$ws.Range("L8").Value = 130000
$ws.Range("N8").Value = -130280
$ws.Range("H40").Value = 4398.5
$ws.Range("I40").Value = 3790
$ws.Range("J40").Value = 4601.3335
$ws.Range("K40").Value = 3790
$ws.Range("L40").Value = 4601.3335
$ws.Range("M40").Value = -3654
$ws.Range("N40").Value = -4873.3335
$ws.Range("H56").Value = 35000
$ws.Range("J56").Value = 35000
$ws.Range("L56").Value = 35000
$ws.Range("N56").Value = -36382
$ws.Range("H80").Value = 58000
$ws.Range("J80").Value = 58000
$ws.Range("L80").Value = 58000
$ws.Range("N80").Value = -60246
$ws.Range("H83").Value = 58000
$ws.Range("J83").Value = 58000
$ws.Range("L83").Value = 174000
$ws.Range("N83").Value = -185232

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 19500
$ws.Range("J49").Value = 19500
$ws.Range("L49").Value = 19500
$ws.Range("N49").Value = -19960
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H116").Value = 233316.33
$ws.Range("J116").Value = 233316.33
$ws.Range("L116").Value = 233316.33
$ws.Range("N116").Value = -242494.33
$ws.Range("H126").Value = 3999.6667
$ws.Range("J126").Value = 3999.5
$ws.Range("L126").Value = 11998.5
$ws.Range("N126").Value = -16938.5
$ws.Range("H132").Value = 3613.2917
$ws.Range("I132").Value = 3462.8096
$ws.Range("K132").Value = 10388.4288
$ws.Range("M132").Value = -7858.4288
$ws.Range("H136").Value = 2740.0278
$ws.Range("I136").Value = 1989.1212
$ws.Range("K136").Value = 5967.363600000001
$ws.Range("M136").Value = -3417.363600000001
